$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying NATMI TPM re-run collapsed the old 6-row FAPs/MuSCs x Lgi3-Stx1a
# pairing table into a 4-row table (rows 2-5) with an additional "Resolving-Mac"
# target cluster, and updated every downstream NATMI statistic. Remove the two
# trailing rows (6 and 7) that no longer exist in the refreshed output.
$ws.Rows("6:7").Delete()

# Row 2: FAPs -> Lgi3/Stx1a -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Lgi3"
$ws.Range("C2").Value = "Stx1a"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.535771
$ws.Range("H2").Value = 4.607313
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.021950333333333
$ws.Range("N2").Value = 9.065851
$ws.Range("O2").Value = 0.4076907863246049
$ws.Range("P2").Value = 0.407690786324605
$ws.Range("Q2").Value = 4.641023685373667
$ws.Range("R2").Value = 41.769213168363
$ws.Range("S2").Value = 0.4076907863246049
$ws.Range("T2").Value = 0.407690786324605

# Row 3: FAPs -> Lgi3/Stx1a -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Lgi3"
$ws.Range("C3").Value = "Stx1a"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.535771
$ws.Range("H3").Value = 4.607313
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.685518
$ws.Range("N3").Value = 8.056554
$ws.Range("O3").Value = 0.3623027595894353
$ws.Range("P3").Value = 0.3623027595894353
$ws.Range("Q3").Value = 4.124340664378001
$ws.Range("R3").Value = 37.119065979402
$ws.Range("S3").Value = 0.3623027595894353
$ws.Range("T3").Value = 0.3623027595894353

# Row 4: FAPs -> Lgi3/Stx1a -> MuSCs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Lgi3"
$ws.Range("C4").Value = "Stx1a"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.535771
$ws.Range("H4").Value = 4.607313
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.692288666666667
$ws.Range("N4").Value = 5.076866
$ws.Range("O4").Value = 0.2283063654591998
$ws.Range("P4").Value = 0.2283063654591998
$ws.Range("Q4").Value = 2.598967857895333
$ws.Range("R4").Value = 23.390710721058
$ws.Range("S4").Value = 0.2283063654591998
$ws.Range("T4").Value = 0.2283063654591998

# Row 5: FAPs -> Lgi3/Stx1a -> Resolving-Mac (new target cluster)
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Lgi3"
$ws.Range("C5").Value = "Stx1a"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.535771
$ws.Range("H5").Value = 4.607313
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01260166666666667
$ws.Range("N5").Value = 0.037805
$ws.Range("O5").Value = 0.001700088626760101
$ws.Range("P5").Value = 0.001700088626760101
$ws.Range("Q5").Value = 0.01935327421833333
$ws.Range("R5").Value = 0.174179467965
$ws.Range("S5").Value = 0.001700088626760101
$ws.Range("T5").Value = 0.001700088626760101
